$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B22").Value = 6289
$ws.Range("C22").Value = 992
$ws.Range("D22").Value = 5839270
$ws.Range("E22").Value = 928.4894259818731
$ws.Range("F22").Value = 8.263040110173868
$ws.Range("G22").Value = 3.765690376569042
$ws.Range("H22").Value = 26.98592541688776
